$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert a new row at position 13, shifting rows 13-23 down to 14-24 ---
# This naturally carries row heights / cell styles the same way Excel does,
# matching the target layout (new dimension A1:C24).
$ws.Rows.Item(13).Insert()

# Copy the cell formatting (styles only) from row 19 (A/B/C all populated, no custom height)
# into the freshly inserted row 13 so B13/C13 get the correct wrap/red styles.
$ws.Range("A19:C19").Copy()
$ws.Range("A13:C13").PasteSpecial(-4122)

# --- Step 2: populate the new row 13 (label-less row: only B13/C13) ---
$ws.Cells.Item(13, 2).Value = "5840793 - Sérgio Schneider"
$ws.Cells.Item(13, 3).Value = "5840793 - Sérgio Schneider"
# Row 13 has no A value in the target -- remove the empty styled cell that Insert() left behind.
$ws.Cells.Item(13, 1).Clear()

# --- Step 3: update cell contents that now hold new/changed text ---

# Row 10 (Objetivos:) - B/C now hold the long objectives paragraph instead of the docente name.
$ws.Cells.Item(10, 2).Value = "1. Descrever os principais processos de conformação mecânica utilizados na indústria metal mecânica.2. Munir o aluno de conhecimentos suficientes para especificar equipamentos e acessórios, usados no processo de conformação, com base nas solicitações mecânicas e variáveis do processo. 3. Ensinar ao aluno a definir tecnicamente o processo adequado de conformação mecânica de produtos da indústria metal mecânica."
$ws.Cells.Item(10, 3).Value = "1. Descrever os principais processos de conformação mecânica utilizados na indústria metal mecânica.2. Munir o aluno de conhecimentos suficientes para especificar equipamentos e acessórios, usados no processo de conformação, com base nas solicitações mecânicas e variáveis do processo. 3. Ensinar ao aluno a definir tecnicamente o processo adequado de conformação mecânica de produtos da indústria metal mecânica."

# Row 14 (Programa resumido:) - B/C now hold the short-syllabus paragraph instead of "Semestral".
$ws.Cells.Item(14, 2).Value = "1. Introdução à Teoria de Plasticidade.2. Classificação dos Processos de Conformação Mecânica.3. Metalurgia da Conformação Mecânica.4. Mecânica da Conformação: Forjamento, Trefilação, Extrusão, Laminação, Dobramento e Estampagem.5. Descrição dos Processos de Conformação Plástica6. Processamento de Chapas Metálicas: Estampagem, Dobramento, Calandragem e corte."
$ws.Cells.Item(14, 3).Value = "1. Introdução à Teoria de Plasticidade.2. Classificação dos Processos de Conformação Mecânica.3. Metalurgia da Conformação Mecânica.4. Mecânica da Conformação: Forjamento, Trefilação, Extrusão, Laminação, Dobramento e Estampagem.5. Descrição dos Processos de Conformação Plástica6. Processamento de Chapas Metálicas: Estampagem, Dobramento, Calandragem e corte."

# Row 16 (Programa:) - B/C now hold the long "Programa" paragraph.
$ws.Cells.Item(16, 2).Value = "1) Introdução à Teoria de Plasticidade: Análise de tensão e deformação 3D, tensões principais, critérios de escoamento, relação tensão-deformação no regime plástico, tensão efetiva e deformação efetiva, energia de deformação e trabalho plástico. 2) Classificação dos processos de conformação. Forjamento, Trefilação, Extrusão Laminação, Estiramento, Estampagem e Dobramento. 3) Metalurgia da Conformação: Noções básicas sobre recozimento, encruamento, conformabilidade, textura e anisotropia. Trabalho a frio e a quente.4) Mecânica da Conformação: métodos de cálculo, efeitos do atrito na conformação e noções básicas de lubrificação. Efeito da taxa de deformação. 5) Descrição dos Processos de Conformação Plástica. Forjamento: Trefilação, Extrusão e Laminação. Equipamentos e acessórios: características e noções de projeto e dimensionamento dos mesmos. Cálculo de carga desses processos.  Análise de defeitos que podem ocorrer nesses processos. Relações geométricas na laminação. Potência e torque de laminação.  6) Processamento de Chapas Metálica. Estampagem: Ensaios de Estampabilidade, Curva Limite de Conformação (CLC). Dobramento: tipos de dobramento, efeito mola, equipamentos e matrizes. Calandragem: tipos de calandragem e equipamentos. Corte: equipamentos de corte e aplicações."
$ws.Cells.Item(16, 3).Value = "1) Introdução à Teoria de Plasticidade: Análise de tensão e deformação 3D, tensões principais, critérios de escoamento, relação tensão-deformação no regime plástico, tensão efetiva e deformação efetiva, energia de deformação e trabalho plástico. 2) Classificação dos processos de conformação. Forjamento, Trefilação, Extrusão Laminação, Estiramento, Estampagem e Dobramento. 3) Metalurgia da Conformação: Noções básicas sobre recozimento, encruamento, conformabilidade, textura e anisotropia. Trabalho a frio e a quente.4) Mecânica da Conformação: métodos de cálculo, efeitos do atrito na conformação e noções básicas de lubrificação. Efeito da taxa de deformação. 5) Descrição dos Processos de Conformação Plástica. Forjamento: Trefilação, Extrusão e Laminação. Equipamentos e acessórios: características e noções de projeto e dimensionamento dos mesmos. Cálculo de carga desses processos.  Análise de defeitos que podem ocorrer nesses processos. Relações geométricas na laminação. Potência e torque de laminação.  6) Processamento de Chapas Metálica. Estampagem: Ensaios de Estampabilidade, Curva Limite de Conformação (CLC). Dobramento: tipos de dobramento, efeito mola, equipamentos e matrizes. Calandragem: tipos de calandragem e equipamentos. Corte: equipamentos de corte e aplicações."

# Row 19 (Método:) - B/C now hold the "Para compor a Nota no Semestre..." text.
$ws.Cells.Item(19, 2).Value = "Para compor a Nota no Semestre (NS) serão feitas duas avaliações (P1 e P2)."
$ws.Cells.Item(19, 3).Value = "Para compor a Nota no Semestre (NS) serão feitas duas avaliações (P1 e P2)."

# Row 20 (Critério:) - B/C now hold the "NS = (P1 + P2)/2..." text.
$ws.Cells.Item(20, 2).Value = "NS = (P1 + P2)/2Serão considerados aprovados os alunos que obtiverem: NS maior ou igual a 5,0. Serão considerados reprovados os alunos que obtiverem: NS menor que 3,0 Para os alunos em que NS é maior ou igual a 3,0 e menor que 5,0 será dada uma prova de recuperação (R)."
$ws.Cells.Item(20, 3).Value = "NS = (P1 + P2)/2Serão considerados aprovados os alunos que obtiverem: NS maior ou igual a 5,0. Serão considerados reprovados os alunos que obtiverem: NS menor que 3,0 Para os alunos em que NS é maior ou igual a 3,0 e menor que 5,0 será dada uma prova de recuperação (R)."

# Row 21 (Norma de recuperação:) - B/C now hold the "A prova de Recuperação..." text.
$ws.Cells.Item(21, 2).Value = "A prova de Recuperação (R) irá compor a nota final (NF) da seguinte forma:NF = (R + NS)/2. Serão considerados aprovados os alunos que obtiverem NF maior ou igual a 5,0."
$ws.Cells.Item(21, 3).Value = "A prova de Recuperação (R) irá compor a nota final (NF) da seguinte forma:NF = (R + NS)/2. Serão considerados aprovados os alunos que obtiverem NF maior ou igual a 5,0."

# Row 22 (Bibliografia:) - B/C now hold the long bibliography text.
$ws.Cells.Item(22, 2).Value = ".  LARKE, E.C. The Rolling of Strip, Sheet, and Plate, Chapman and Hall, 19672.  HONEYCOMBE, R.W.K. The Plastic Deformation of Metals, Edward Arnold, 1968.3.  HOSFORD, W.F. Metal Forming: Mechanics and Metallurgy, Prentice-Hall, 1983. 4.  WEERTMAN, J. Elementary Dislocation Theory, Collier-McMillan, 1965. 5.  AVITZUR, B. Metal Forming: Processes and Analysis, McGraw-Hill, 1968.6. BRESCIANI Filho, E. e outros. Conformação Plástica dos Metais, Editora da UNICAMP           Campinas, Volumes 1 e 2, 1986. 7. CETLIN, P. R.; HELMAN, H. Fundamentos de Conformação Mecânica dos Metais. Art Liber: São Paulo, 2005.8. ROWE, G.W. Elements of Metalworking Theory. Edward Arnold Publishers, 1979.9.  JOHNSON, W.; MELLOR, P.B. Engineering Plasticity, Van Nostrand Reinhold, 1973. 10. DIETER, G. E. Metalurgia Mecânica. Guanabara Dois, 1981.11. SCHAEFFER, L. Introdução à Conformação Mecânica dos Metais, Ed. da UFRGS, 1983. 12. RODRIGUES, J. Tecnologia Mecânica. Volumes 1 e 2, Ed. Escolar, 2005. 13. CALLISTER, W. D. Ciência e Engenharia dos Materiais: Uma Introdução. Rio de Janeiro: LTC, 1999."
$ws.Cells.Item(22, 3).Value = ".  LARKE, E.C. The Rolling of Strip, Sheet, and Plate, Chapman and Hall, 19672.  HONEYCOMBE, R.W.K. The Plastic Deformation of Metals, Edward Arnold, 1968.3.  HOSFORD, W.F. Metal Forming: Mechanics and Metallurgy, Prentice-Hall, 1983. 4.  WEERTMAN, J. Elementary Dislocation Theory, Collier-McMillan, 1965. 5.  AVITZUR, B. Metal Forming: Processes and Analysis, McGraw-Hill, 1968.6. BRESCIANI Filho, E. e outros. Conformação Plástica dos Metais, Editora da UNICAMP           Campinas, Volumes 1 e 2, 1986. 7. CETLIN, P. R.; HELMAN, H. Fundamentos de Conformação Mecânica dos Metais. Art Liber: São Paulo, 2005.8. ROWE, G.W. Elements of Metalworking Theory. Edward Arnold Publishers, 1979.9.  JOHNSON, W.; MELLOR, P.B. Engineering Plasticity, Van Nostrand Reinhold, 1973. 10. DIETER, G. E. Metalurgia Mecânica. Guanabara Dois, 1981.11. SCHAEFFER, L. Introdução à Conformação Mecânica dos Metais, Ed. da UFRGS, 1983. 12. RODRIGUES, J. Tecnologia Mecânica. Volumes 1 e 2, Ed. Escolar, 2005. 13. CALLISTER, W. D. Ciência e Engenharia dos Materiais: Uma Introdução. Rio de Janeiro: LTC, 1999."
